# Daily TGP (terminal gate pricing) roll-forward.
#
# Each state block on Sheet1 lists N terminals for "today" followed by the
# same N terminals for "yesterday" (2*N rows per block). A new day's prices
# arrive: the old "today" rows become "yesterday" (shifted down by N rows,
# values unchanged) and the old "yesterday" rows are discarded; brand new
# prices for the new "today" (effective date 46009, i.e. one day after the
# previous top date of 46008) are written into the freed top N rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each block: starting row of the "today" rows, and how many terminal rows
# (N) it contains (so the "yesterday" rows run from start+N .. start+2N-1).
$blocks = @(
    @{ Start = 8;  N = 3 },   # New South Wales
    @{ Start = 17; N = 1 },   # Northern Territory
    @{ Start = 22; N = 5 },   # Queensland
    @{ Start = 35; N = 1 },   # South Australia
    @{ Start = 40; N = 2 },   # Tasmania
    @{ Start = 47; N = 2 },   # Victoria
    @{ Start = 54; N = 6 }    # Western Australia
)

# New "today" data for each top row, keyed by row number. Missing D/E/F/G
# keys mean "leave that cell alone" (it's a non-numeric N/A cell in the
# original sheet and stays N/A).
$newData = @{
    8  = @{ D = 156.56; E = 156.88999999999999; F = 166.89; G = 157 }
    9  = @{ D = 156.56; E = 156.88999999999999; F = 166.89; G = 157 }
    10 = @{ D = 158.63; E = 159.46;             F = 169.46; G = 159.97999999999999 }

    17 = @{ D = 161.69; E = 163.15; F = 173.15 }

    22 = @{ D = 157.37;             E = 158.46; F = 168.06; G = 159.62 }
    23 = @{ D = 163.55000000000001; E = 163.57; F = 173.57 }
    24 = @{ D = 163.32;             E = 163.89; F = 173.89 }
    25 = @{ D = 163.82;             E = 163.53; F = 173.53; G = 163.30000000000001 }
    26 = @{ D = 162.72999999999999; E = 164.96; F = 174.96 }

    35 = @{ D = 156.91999999999999; E = 156.6; F = 165.6 }

    40 = @{ D = 163.02000000000001; E = 163.81; F = 173.81 }
    41 = @{ D = 162.74;             E = 164.23; F = 174.23 }

    47 = @{ D = 158.11000000000001; E = 158.53; F = 168.53 }
    48 = @{ D = 157.91999999999999; E = 158.62; F = 168.62 }

    54 = @{ D = 172.51; E = 174.33; F = 184.33 }
    55 = @{ D = 160.71; E = 169.25; F = 179.25 }
    56 = @{ D = 162.97 }
    57 = @{ D = 161.99; E = 163.52000000000001 }
    58 = @{ D = 157.88999999999999; E = 159.57; F = 169.57 }
    59 = @{ D = 164.61; E = 171.59 }
}

# The new effective date (one day later than the previous top date).
$newDate = 46009

foreach ($block in $blocks) {
    $start = $block.Start
    $n = $block.N

    # 1) Shift the current "today" rows down into the "yesterday" rows
    #    (copy Effective Date + Diesel/ULP/PULP/e10 values as-is; terminal
    #    name in column C already matches row-for-row between the two
    #    halves of the block, so it is left untouched).
    for ($i = 0; $i -lt $n; $i++) {
        $srcRow = $start + $i
        $dstRow = $start + $n + $i
        foreach ($col in @("A", "D", "E", "F", "G")) {
            $srcCell = $ws.Range("$col$srcRow")
            $srcVal = $srcCell.Value2
            if ($srcVal -ne $null) {
                $ws.Range("$col$dstRow").Value = $srcVal
            }
        }
    }

    # 2) Write the brand new "today" data into the freed top rows.
    for ($i = 0; $i -lt $n; $i++) {
        $row = $start + $i
        $ws.Range("A$row").Value = $newDate
        $vals = $newData[$row]
        foreach ($col in @("D", "E", "F", "G")) {
            if ($vals.ContainsKey($col)) {
                $ws.Range("$col$row").Value = $vals[$col]
            }
        }
    }
}
